# 14/07 - Commit EOD
# Add two new "SOI1312" test-company rows into GeneralVariables and make
# that sheet the active one (mirrors a user editing + saving that tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneralVariables")

# New company block: insert the "name" pair right after the other
# testingCompanySOI* / AutoTestingCompany_SOI* rows (old row 16 -> new row 16)
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "testingCompanySOI1312"
$ws.Range("B16").Value = "AutoTestingCompany_SOI1312"

# New company block: insert the matching "id" pair right after the other
# idTestingCompanySOI* rows (old row 30 -> new row 31, since row 16 already
# shifted everything below it down by one)
$ws.Rows.Item(31).Insert()
$ws.Range("B31").Value = "0013E00001CtYgQQAV"
$ws.Range("A31").Value = "idTestingCompanySOI1312"

# Make GeneralVariables the active/selected sheet with B16 as the active
# cell, matching the saved selection state.
$ws.Activate()
$ws.Range("B16").Select()
